$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weekly AC/AN shared-string header labels (16/08 -> 23/08)
$ws.Range("F1").Value = "AC_20200823"
$ws.Range("G1").Value = "AN_20200823"

# Refreshed monitoring values for columns C:G (AN_ANUAL, AN_SEM, AN_TRI, AC, AN), rows 2-42
$data = @(
    @(-8.3, -4.9, 30.7, 130.7, 42.3),
    @(-4.2, 6.7, 52.1, 145.2, 42.7),
    @(-7.1, 1.7, 18.6, 132.1, 109.4),
    @(-24.9, -11, 22.5, 142.3, 46.7),
    @(-23.1, -22.9, -9.1, 136.6, 41.9),
    @(-21.2, -22.2, -11.5, 137.6, 40.8),
    @(-49.3, -37.2, 5.8, 151.1, 79.4),
    @(-15.9, -2.3, 38.2, 123.1, 11),
    @(-16.8, -18.2, 2.2, 155.9, 62),
    @(-10.5, -11.6, 6.2, 132.3, 72),
    @(-11.1, -12.1, 8.6, 132.6, 71.2),
    @(-10.5, -11.1, 18.7, 129.5, 51.5),
    @(-19.8, -16.6, 23.2, 108.8, 12.8),
    @(-25, -21.5, 26, 138.1, 122.6),
    @(-2.6, 10.6, 56.7, 146.6, 37),
    @(-15.1, -4, 35.7, 131.5, 22.2),
    @(-26.4, -26.5, -8.2, 156.6, 108.2),
    @(-14.8, -14.9, 1.1, 133.6, 56.3),
    @(-41.3, -47.5, -24.5, 130.1, 53.1),
    @(-17.9, -10.9, 30.6, 175.2, 183.7),
    @(-33.6, -26.3, 6.5, 178.6, 106.3),
    @(-24.6, -21.9, -10.9, 146.4, 62.2),
    @(-50, -46, -10.6, 107, 15.3),
    @(-52.1, -47.3, -12.1, 111.7, 10.7),
    @(-45.5, -37.7, 0.7, 101.5, 14),
    @(-32.5, -24.2, 20.3, 206.5, 143.8),
    @(-24.9, -9.3, 42.4, 199.3, 99.6),
    @(-30.8, -25.7, 14.2, 161.5, 105),
    @(-21.4, -23.8, -13.4, 140.7, 57.9),
    @(-23.5, -28.7, -16.7, 156.6, 60.6),
    @(-46.5, -50.1, -20.6, 116.8, 51),
    @(-26.6, -31.5, 15.7, 121.1, 175.9),
    @(-44.6, -47.6, -6.7, 90.9, -7.5),
    @(-48, -50.3, -8.4, 91.6, 32.4),
    @(-43.3, -45.8, -13, 106.2, 22.1),
    @(-43.1, -44.8, -12.4, 104.8, 23.7),
    @(-48.6, -51.3, -10.3, 94, 30),
    @(-46.5, -46, -10.9, 102.2, 26.9),
    @(-57.5, -56.6, -18.5, 71.3, -31.5),
    @(-38.3, -43.2, -8.7, 129.8, 64.4),
    @(-30.6, -27.7, -8.5, 96.8, 12.1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item(2 + $r, 3 + $c).Value = $rowVals[$c]
    }
}
